{"js": "// Update \"Lista de control\" document for the new edition of the contest:\n// - bump edition number XIII -> XIV\n// - bump registration deadline year 2023 -> 2024\n// - bump final date from \"24 de mayo de 2024\" to \"lunes 19 de mayo de 2025\"\n//   (first two mentions get the weekday \"lunes\"; the third mention of the\n//   date only updates the day/year)\n// - bump \"01 de marzo de 2024\" -> \"01 de marzo de 2025\"\n// - bump \"30 de abril de 2024\" -> \"30 de abril de 2025\"\n// - remove \"tanto en una semifinal como\" from the prize-eligibility rule\n\nconst body = context.document.body;\n\n// 1) Edition number: XIII -> XIV\nlet res = body.search(\"XIII edici\u00f3n\", { matchCase: true });\nres.load(\"text\");\nawait context.sync();\nif (res.items.length > 0) {\n  res.items[0].insertText(\"XIV edici\u00f3n\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 2) Registration deadline year: 2023 -> 2024\nres = body.search(\"30 de noviembre de 2023\", { matchCase: true });\nres.load(\"text\");\nawait context.sync();\nif (res.items.length > 0) {\n  res.items[0].insertText(\"30 de noviembre de 2024\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 3) Final date \"24 de mayo de 2024\" appears 3 times, in document order:\n//    a) \"... Instituto Cervantes el 24 de mayo de 2024\"        -> add \"lunes \"\n//    b) \"Antes de la final en Fr\u00e1ncfort el 24 de mayo de 2024\" -> add \"lunes \"\n//    c) \"... que el d\u00eda 24 de mayo de 2024 estar\u00e9is en ...\"    -> no \"lunes\"\n// All three also bump the year 2024 -> 2025.\nres = body.search(\"24 de mayo de 2024\", { matchCase: true });\nres.load(\"text\");\nawait context.sync();\n\nconst replacements = [\n  \"lunes 19 de mayo de 2025\",\n  \"lunes 19 de mayo de 2025\",\n  \"19 de mayo de 2025\",\n];\nfor (let i = 0; i < res.items.length && i < replacements.length; i++) {\n  res.items[i].insertText(replacements[i], Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 4) Selection date for semifinal candidates: 2024 -> 2025\nres = body.search(\"01 de marzo de 2024\", { matchCase: true });\nres.load(\"text\");\nawait context.sync();\nif (res.items.length > 0) {\n  res.items[0].insertText(\"01 de marzo de 2025\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 5) Semifinal deadline: 2024 -> 2025\nres = body.search(\"30 de abril de 2024\", { matchCase: true });\nres.load(\"text\");\nawait context.sync();\nif (res.items.length > 0) {\n  res.items[0].insertText(\"30 de abril de 2025\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 6) Prize-eligibility rule: drop \"tanto en una semifinal como\"\nres = body.search(\"tanto en una semifinal como en la final.\", { matchCase: true });\nres.load(\"text\");\nawait context.sync();\nif (res.items.length > 0) {\n  res.items[0].insertText(\"en la final.\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Update \"Lista de control\" document for the new edition of the contest:\n# - bump edition number XIII -> XIV\n# - bump registration deadline year 2023 -> 2024\n# - bump final date from \"24 de mayo de 2024\" to \"lunes 19 de mayo de 2025\"\n#   (first two mentions get the weekday \"lunes\"; the third mention of the\n#   date only updates the day/year)\n# - bump \"01 de marzo de 2024\" -> \"01 de marzo de 2025\"\n# - bump \"30 de abril de 2024\" -> \"30 de abril de 2025\"\n# - remove \"tanto en una semifinal como\" from the prize-eligibility rule\n\n$d = $word.ActiveDocument\n\nfunction Replace-NthOccurrence {\n    param(\n        $doc,\n        [string]$SearchText,\n        [int]$Occurrence,\n        [string]$ReplaceText\n    )\n    $rng = $doc.Content\n    $rng.Start = 0\n    $rng.End = $doc.Content.End\n    $count = 0\n    while ($true) {\n        $found = $rng.Find.Execute($SearchText, $false, $false, $false, $false, $false, $true, 1, $false, $null, 0)\n        if (-not $found) { break }\n        $count = $count + 1\n        if ($count -eq $Occurrence) {\n            $rng.Text = $ReplaceText\n            return $true\n        }\n        $rng.Start = $rng.End\n        $rng.End = $doc.Content.End\n    }\n    return $false\n}\n\n# 1) Edition number: XIII -> XIV\nReplace-NthOccurrence $d \"XIII edici\u00f3n\" 1 \"XIV edici\u00f3n\" | Out-Null\n\n# 2) Registration deadline year: 2023 -> 2024\nReplace-NthOccurrence $d \"30 de noviembre de 2023\" 1 \"30 de noviembre de 2024\" | Out-Null\n\n# 3) Final date \"24 de mayo de 2024\" appears 3 times, in document order:\n#    a) \"... Instituto Cervantes el 24 de mayo de 2024\"        -> add \"lunes \"\n#    b) \"Antes de la final en Fr\u00e1ncfort el 24 de mayo de 2024\" -> add \"lunes \"\n#    c) \"... que el d\u00eda 24 de mayo de 2024 estar\u00e9is en ...\"    -> no \"lunes\"\n# All three also bump the year 2024 -> 2025. Replace occurrence #1 each time\n# since earlier matches disappear once replaced.\nReplace-NthOccurrence $d \"24 de mayo de 2024\" 1 \"lunes 19 de mayo de 2025\" | Out-Null\nReplace-NthOccurrence $d \"24 de mayo de 2024\" 1 \"lunes 19 de mayo de 2025\" | Out-Null\nReplace-NthOccurrence $d \"24 de mayo de 2024\" 1 \"19 de mayo de 2025\" | Out-Null\n\n# 4) Selection date for semifinal candidates: 2024 -> 2025\nReplace-NthOccurrence $d \"01 de marzo de 2024\" 1 \"01 de marzo de 2025\" | Out-Null\n\n# 5) Semifinal deadline: 2024 -> 2025\nReplace-NthOccurrence $d \"30 de abril de 2024\" 1 \"30 de abril de 2025\" | Out-Null\n\n# 6) Prize-eligibility rule: drop \"tanto en una semifinal como\"\nReplace-NthOccurrence $d \"tanto en una semifinal como en la final.\" 1 \"en la final.\" | Out-Null\n"}
